$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the value text for the two new rows first (matches original authoring order)
$ws.Range("B14").Value = "ATTACK"
$ws.Range("B15").Value = "DEFEND"

# Update the title value (row 4, column B) with the new text
$ws.Range("B4").Value = "Pengu and the Unlikely Encounters"

# Add the key text for the two new rows: attack/defend
$ws.Range("A14").Value = "attack"
$ws.Range("A15").Value = "defend"

# Update the selected cell to reflect the new active cell after edits
$ws.Range("A16").Select()
